$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Simplify the Treatment Agent query: CONCAT() around a single argument was
#    redundant, so drop the wrapper and keep the REPLACE() call.
$old = "CONCAT(REPLACE(trt.treatment_agent, ';', ', ')) AS ""Treatment Agent"","
$new = "REPLACE(trt.treatment_agent, ';', ', ') AS ""Treatment Agent"","
$cell = $ws.Range("B5")
$cell.Value2 = $cell.Value2.Replace($old, $new)

# 2. The Studies/Participants/Diagnosis query cells (B2:B4) get re-formatted
#    to match the Treatment query cell (B5): same wrapped, size-12 Calibri text.
$ws.Range("B5").Copy()
$ws.Range("B2:B4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# 3. Move the selection to C5 (the view had scrolled down to row 4).
$ws.Range("C5").Select()
